$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.164.32"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "2.501.72"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.07"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -3.86%  "

$ws.Range("E10").Value = "  -3.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.99"
$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("D15").Value = "2.892.88"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "2.504.50"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("E17").Value = "  -2.42%  "

$ws.Range("D18").Value = "48.041.55"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.97"
$ws.Range("E20").Value = "  +8.48%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.23"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.62"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -1.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.84"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  +4.45%  "

$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("E30").Value = "  -4.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.08"
$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.32"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("E34").Value = "  -4.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0778"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("E37").Value = "  -2.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.59"
$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.52"
$ws.Range("E41").Value = "  +2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.01"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0306"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("D45").Value = "2.009.11"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").Value = "  +2.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.95"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.96"
$ws.Range("E51").Value = "  -1.31%  "
